$d = $word.ActiveDocument

$pkgPre = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgPost = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParagraphContainingText($doc, $searchText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return $null
    }
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphIndexContainingText($doc, $searchText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return -1
    }
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $i
        }
    }
    return -1
}

function Replace-ParagraphXml($doc, $para, $bodyXml) {
    $xml = $pkgPre + $bodyXml + $pkgPost
    $para.Range.InsertXML($xml)
}

# 1) "Explore and format site data..." -> "Explore and format **site** data..."
$p1 = Get-ParagraphContainingText $d "Explore and format site data"
Replace-ParagraphXml $d $p1 '<w:p w14:paraId="230287ED" w14:textId="77777777" w:rsidR="005E2611" w:rsidRDefault="005E2611" w:rsidP="005E2611"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Explore and format </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>site</w:t></w:r><w:r><w:t xml:space="preserve"> data. Sites should be your first step in data formatting. Prior to following the formatting steps in Section one of this document, take a moment to explore how sites are coded. Of importance are:</w:t></w:r></w:p>
'

# 2) Large paragraph: ddply sentence insertion + bookmark relocation + run merge ("Again, i" + "f the problem")
$p2 = Get-ParagraphContainingText $d "How many records are there per site?"
Replace-ParagraphXml $d $p2 '<w:p w14:paraId="4E2E32CE" w14:textId="4605D04D" w:rsidR="004C17AB" w:rsidRDefault="00155338" w:rsidP="004C17AB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>How many records are there per site?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C645D9"><w:t xml:space="preserve">Sites that are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C645D9"><w:t>mis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00C645D9"><w:t xml:space="preserve">-defined can also be determined by observing the number of records across sites. </w:t></w:r><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">If sites are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>mis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">-defined, this can be identified if a large proportion of sites have very few records. </w:t></w:r><w:r w:rsidR="00C645D9"><w:t xml:space="preserve">There are many ways to </w:t></w:r><w:r w:rsidR="004C17AB"><w:t>determine</w:t></w:r><w:r w:rsidR="00C645D9"><w:t xml:space="preserve"> this. To observe the number of records per site using the table function in base R, </w:t></w:r><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">use either </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="004C17AB"><w:t>table(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004C17AB"><w:t>example_df$site</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">) to observe the records in wide format or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>data.frame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t>(table(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>example_df$site</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">)) to observe the records in long format. </w:t></w:r><w:r><w:t xml:space="preserve">The latter can also be done in Hadley Wickham’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plyr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package using: </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>ddply</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>example_df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,.(site</w:t></w:r><w:r><w:t>),''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nrow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>'')</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">If there are a large number of sites, it can be cumbersome to search through them all. You can avoid this by ordering from the smallest to largest number of records per site. First, assign a name to your site table: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>xy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve"> &lt;- </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidR="004C17AB"><w:t>dat</w:t></w:r><w:r w:rsidR="004C17AB"><w:t>a.frame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004C17AB"><w:t>table(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>example_df$site</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">)). Next, order by frequency: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="004C17AB" w:rsidRPr="004C17AB"><w:t>xy2[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004C17AB" w:rsidRPr="004C17AB"><w:t>order(xy2$Freq),]</w:t></w:r><w:r w:rsidR="004C17AB"><w:t>. As above, modify as ne</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">cessary, provide descriptive comments in your script for your modification, and add-commit-push to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>GitHub</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="004C17AB"><w:lastRenderedPageBreak/><w:t xml:space="preserve">Again, if the problem is not clear to you, add an issue to the core-transient </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C17AB"><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve"> hub </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="004C17AB"><w:t>repository,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004C17AB"><w:t xml:space="preserve"> describe the problem in detail and assign the issue to me.</w:t></w:r></w:p>'

# 3) "Once you are " / "done with site exploration" / ", save your script and " -> merged into a single run
$idx3 = Get-ParagraphIndexContainingText $d "done with site exploration"
$p3 = $d.Paragraphs.Item($idx3)
Replace-ParagraphXml $d $p3 '<w:p w14:paraId="0470400B" w14:textId="0E3A9BB5" w:rsidR="004C17AB" w:rsidRDefault="004C17AB" w:rsidP="004C17AB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Once you are done with site exploration, save your script and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-add-commit-push.</w:t></w:r></w:p>
'

# 4) Empty paragraph (two after paragraph 3) gets new text: "Explore and format **species** data."
$p4 = $d.Paragraphs.Item($idx3 + 2)
Replace-ParagraphXml $d $p4 '<w:p w14:paraId="7809FAD3" w14:textId="77777777" w:rsidR="004C17AB" w:rsidRDefault="004C17AB" w:rsidP="004C17AB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Explore and format </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>species</w:t></w:r><w:r><w:t xml:space="preserve"> data.</w:t></w:r></w:p>
'

# 5) Final paragraph (three after paragraph 3) loses the _GoBack bookmark (it moved into paragraph 2 above)
$p5 = $d.Paragraphs.Item($idx3 + 3)
Replace-ParagraphXml $d $p5 '<w:p w14:paraId="4EC08807" w14:textId="2F9DA776" w:rsidR="00155338" w:rsidRPr="00386BF3" w:rsidRDefault="00155338" w:rsidP="004C17AB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p>
'

